$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.560.44'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.754.95'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.16%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '324.94'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('E6').Value = '  -0.14%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4572'
$c.ClearFormats()
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('E8').Value = '  -1.78%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07449'
$c.ClearFormats()
$ws.Range('E9').Value = '  -0.52%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '41.52'
$c.ClearFormats()
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('E12').Value = '  -0.10%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '20.78'
$c.ClearFormats()
$ws.Range('E13').Value = '  +0.56%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.006'
$c.ClearFormats()
$ws.Range('E14').Value = '  -0.56%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.160'
$c.ClearFormats()
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').Value = '1.754.47'
$ws.Range('E16').Value = '  +0.19%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '93.49'
$c.ClearFormats()
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('E18').Value = '  -0.82%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06409'
$c.ClearFormats()
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('D23').Value = '27.597.07'
$ws.Range('E23').Value = '  +0.01%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.21'
$c.ClearFormats()
$ws.Range('E25').Value = '  -1.32%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '165.30'
$c.ClearFormats()
$ws.Range('E26').Value = '  +2.16%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '20.11'
$c.ClearFormats()
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('D28').Value = '1.952.61'
$ws.Range('E28').Value = '  +0.01%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.133'
$c.ClearFormats()
$ws.Range('E29').Value = '  +0.83%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '125.34'
$c.ClearFormats()
$ws.Range('E30').Value = '  +0.07%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.075'
$c.ClearFormats()
$ws.Range('E31').Value = '  -0.57%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.09217'
$c.ClearFormats()
$ws.Range('E32').Value = '  +2.15%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.663'
$c.ClearFormats()
$ws.Range('E33').Value = '  -0.31%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.505'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.57%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '11.73'
$c.ClearFormats()
$ws.Range('E35').Value = '  -2.15%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.02273'
$c.ClearFormats()
$ws.Range('E36').Value = '  -1.93%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.06020'
$c.ClearFormats()
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.2083'
$c.ClearFormats()
$ws.Range('E38').Value = '  -0.33%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.6263'
$c.ClearFormats()
$ws.Range('E39').Value = '  -1.45%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '4.921'
$c.ClearFormats()
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('E42').Value = '  -0.95%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '7.756'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.22%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '13.12'
$c.ClearFormats()
$ws.Range('E44').Value = '  -1.58%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.716'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.07%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.5861'
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.48%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '121.86'
$c.ClearFormats()
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('E48').Value = '  -0.90%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.06899'
$c.ClearFormats()
$ws.Range('E49').Value = '  +0.33%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.130'
$c.ClearFormats()
$ws.Range('E50').Value = '  -2.41%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '72.03'
$c.ClearFormats()
$ws.Range('E51').Value = '  -0.32%  '
